$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price/volume snapshot (GitHub Actions scheduled update).
# Price cells that would otherwise be auto-parsed as numbers (and lose trailing
# zeros) are entered with a leading apostrophe so Excel stores them as text,
# matching how the sheet already stores every Price/Volume value as a string.

$ws.Range("D2").Value = '59.959.73'
$ws.Range("E2").Value = '  -0.81%  '
$ws.Range("D3").Value = '2.657.05'
$ws.Range("E3").Value = '  +1.20%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '''522.16'
$ws.Range("E5").Value = '  -0.27%  '
$ws.Range("D6").Value = '''148.85'
$ws.Range("E6").Value = '  -0.93%  '
$ws.Range("D7").Value = '''0.995'
$ws.Range("E7").Value = '  -0.38%  '
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("D9").Value = '2.691.18'
$ws.Range("E9").Value = '  +2.15%  '
$ws.Range("E10").Value = '  +2.73%  '
$ws.Range("D11").Value = '''0.106'
$ws.Range("E11").Value = '  +0.57%  '
$ws.Range("E12").Value = '  -0.56%  '
$ws.Range("E13").Value = '  -1.38%  '
$ws.Range("D14").Value = '3.126.86'
$ws.Range("E14").Value = '  +1.10%  '
$ws.Range("D15").Value = '59.910.06'
$ws.Range("E15").Value = '  -0.91%  '
$ws.Range("D16").Value = '''21.57'
$ws.Range("E16").Value = '  +0.26%  '
$ws.Range("E17").Value = '  +0.65%  '
$ws.Range("D18").Value = '2.678.95'
$ws.Range("E18").Value = '  +1.64%  '
$ws.Range("D19").Value = '''4.64'
$ws.Range("E19").Value = '  -0.22%  '
$ws.Range("D20").Value = '''349.13'
$ws.Range("E20").Value = '  +0.84%  '
$ws.Range("D21").Value = '''10.66'
$ws.Range("E21").Value = '  +1.70%  '
$ws.Range("D22").Value = '''6.27'
$ws.Range("E22").Value = '  +1.65%  '
$ws.Range("D23").Value = '''0.997'
$ws.Range("E23").Value = '  +0.30%  '
$ws.Range("D24").Value = '''61.54'
$ws.Range("E24").Value = '  +1.05%  '
$ws.Range("D25").Value = '''0.429'
$ws.Range("E25").Value = '  +1.87%  '
$ws.Range("D26").Value = '2.764.54'
$ws.Range("E26").Value = '  +0.45%  '
$ws.Range("B27").Value = 'Binance-PegBSC-USD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D27").Value = '''0.997'
$ws.Range("E27").Value = '  -0.18%  '
$ws.Range("B28").Value = 'Kaspa'
$ws.Range("C28").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D28").Value = '''0.162'
$ws.Range("E28").Value = '  -1.07%  '
$ws.Range("D29").Value = '0.0₃0837'
$ws.Range("E29").Value = '  +1.37%  '
$ws.Range("D30").Value = '''7.22'
$ws.Range("E30").Value = '  +1.68%  '
$ws.Range("D31").Value = '''6.62'
$ws.Range("E31").Value = '  +10.27%  '
$ws.Range("D32").Value = '''0.997'
$ws.Range("E32").Value = '  -0.30%  '
$ws.Range("D33").Value = '''1.60'
$ws.Range("E33").Value = '  -0.05%  '
$ws.Range("D34").Value = '''19.13'
$ws.Range("E34").Value = '  +0.40%  '
$ws.Range("E35").Value = '  +18.87%  '
$ws.Range("D36").Value = '''149.76'
$ws.Range("E36").Value = '  -0.67%  '
$ws.Range("E37").Value = '  +2.12%  '
$ws.Range("E38").Value = '  +0.92%  '
$ws.Range("D39").Value = '''0.882'
$ws.Range("E39").Value = '  +0.70%  '
$ws.Range("D40").Value = '''36.75'
$ws.Range("E40").Value = '  +0.34%  '
$ws.Range("D41").Value = '''3.73'
$ws.Range("E41").Value = '  +1.82%  '
$ws.Range("E42").Value = '  -0.23%  '
$ws.Range("D43").Value = '''290.43'
$ws.Range("E43").Value = '  -0.21%  '
$ws.Range("D44").Value = '''0.630'
$ws.Range("E44").Value = '  -0.39%  '
$ws.Range("D45").Value = '''0.1000'
$ws.Range("E45").Value = '  -0.55%  '
$ws.Range("D46").Value = '''0.994'
$ws.Range("E46").Value = '  -0.50%  '
$ws.Range("D47").Value = '''19.81'
$ws.Range("E47").Value = '  -0.26%  '
$ws.Range("D48").Value = '''0.0551'
$ws.Range("E48").Value = '  -0.04%  '
$ws.Range("E49").Value = '  +1.16%  '
$ws.Range("D50").Value = '''0.0234'
$ws.Range("E50").Value = '  -0.38%  '
$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").Value = '1.996.65'
$ws.Range("E51").Value = '  +1.54%  '
